$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.454.47"
$ws.Range("E2").Value = "  +0.13%  "
$ws.Range("D3").Value = "2.072.91"
$ws.Range("E3").Value = "  +0.55%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.24"
$ws.Range("E5").Value = "  -0.77%  "
$ws.Range("E6").Value = "  +1.66%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "57.44"
$ws.Range("E8").Value = "  -1.33%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.397"
$ws.Range("E9").Value = "  +3.53%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0775"
$ws.Range("E10").Value = "  +1.76%  "
$ws.Range("E11").Value = "  +0.92%  "
$ws.Range("D12").Value = "2.376.80"
$ws.Range("E12").Value = "  +0.53%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.49"
$ws.Range("E13").Value = "  +0.61%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.61"
$ws.Range("E14").Value = "  -2.73%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.780"
$ws.Range("E15").Value = "  -0.04%  "
$ws.Range("E16").Value = "  -0.01%  "
$ws.Range("D17").Value = "2.072.27"
$ws.Range("E17").Value = "  +0.51%  "
$ws.Range("D18").Value = "37.384.67"
$ws.Range("E18").Value = "  -0.53%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.42"
$ws.Range("E19").Value = "  +4.23%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.84"
$ws.Range("E20").Value = "  +1.17%  "
$ws.Range("E21").Value = "  +0.64%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "227.31"
$ws.Range("E22").Value = "  +0.79%  "
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("E24").Value = "  +0.99%  "
$ws.Range("E25").Value = "  -2.49%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "166.88"
$ws.Range("E26").Value = "  +1.89%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.85"
$ws.Range("E27").Value = "  -0.42%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.45"
$ws.Range("E28").Value = "  -0.61%  "
$ws.Range("E29").Value = "  +0.91%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.16"
$ws.Range("E30").Value = "  -0.23%  "
$ws.Range("E31").Value = "  -0.39%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.55"
$ws.Range("E32").Value = "  +0.94%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0621"
$ws.Range("E34").Value = "  +2.29%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.50"
$ws.Range("E35").Value = "  -3.63%  "
$ws.Range("E36").Value = "  +0.73%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.32"
$ws.Range("E38").Value = "  +0.00%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.72"
$ws.Range("E39").Value = "  -2.66%  "
$ws.Range("E40").Value = "  -0.69%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "98.97"
$ws.Range("E41").Value = "  +1.41%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0962"
$ws.Range("E42").Value = "  -2.18%  "
$ws.Range("D43").Value = "1.482.95"
$ws.Range("E43").Value = "  +0.08%  "
$ws.Range("E44").Value = "  +2.49%  "
$ws.Range("E45").Value = "  +1.41%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.10"
$ws.Range("E46").Value = "  -8.92%  "
$ws.Range("E47").Value = "  +0.46%  "
$ws.Range("E48").Value = "  -3.90%  "
$ws.Range("E49").Value = "  +0.64%  "
$ws.Range("E50").Value = "  +0.60%  "
$ws.Range("D51").Value = "2.263.34"
$ws.Range("E51").Value = "  +0.48%  "
